$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '(''Angel'', [''Token Creature — Angel'', ''Flying'', ''4/4''])'
$ws.Range("A3").Value = '(''Beast'', [''Token Creature — Beast'', ''3/3''])'
$ws.Range("A4").Value = '(''Cat'', [''Token Creature — Cat'', ''2/2''])'
$ws.Range("A5").Value = '(''Dragon'', [''Token Creature — Dragon'', ''Flying'', ''{R}: This creature gets +1/+0 until end of turn.'', ''2/2''])'
$ws.Range("A6").Value = '(''Elemental'', [''Token Creature — Elemental'', ''1/1''])'
$ws.Range("A7").Value = '(''Garruk, Caller of Beasts Emblem'', [''Emblem — Garruk'', ''Whenever you cast a creature spell, you may search your library for a creature card, put it onto the battlefield, then shuffle your library.''])'
$ws.Range("A8").Value = '(''Goat'', [''Token Creature — Goat'', ''0/1''])'
$ws.Range("A9").Value = '(''Liliana of the Dark Realms Emblem'', [''Emblem — Liliana'', "Swamps you control have ‘{T}: Add {B}{B}{B}{B}.''"])'
$ws.Range("A10").Value = '(''Saproling'', [''Token Creature — Saproling'', ''1/1''])'
$ws.Range("A11").Value = '(''Sliver'', [''Token Creature — Sliver'', ''1/1''])'
$ws.Range("A12").Value = '(''Wolf'', [''Token Creature — Wolf'', ''2/2''])'
$ws.Range("A13").Value = '(''Zombie'', [''Token Creature — Zombie'', ''2/2''])'

# Remove now-empty trailing rows 14:43 (previous data consolidated into A2:A13)
$ws.Range("A14:A43").EntireRow.Delete() | Out-Null
